$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value added
$ws1.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact" row (old row 11) which shifts subsequent rows up
$ws1.Rows.Item(11).Delete()

# Row 10 (was "Contact"/"No display for ContactDetail") becomes "Jurisdiction"/"United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Case Sensitive value added
$ws1.Range("B14").Value = "true"

# --- Concepts sheet (no content changes, only shared string bookkeeping) ---
